$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Rebalance creature HP (column E) for rows 12-16 from 70 to 50
$ws.Range("E12:E16").Value = 50

# Update the active selection to F12, matching the saved selection state
$ws.Range("F12").Select()
